# Updates the cryptocurrency price/volume table on the active sheet
# to reflect the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" column values are stored as TEXT (they look like
# numbers, e.g. "0.9972" or "1.208", but the source sheet keeps them as
# inline strings). A leading apostrophe forces Excel to keep them as text
# instead of silently converting them to numeric values.

# Row 2
$ws.Range("D2").Value = '25.779.90'
$ws.Range("E2").Value = '  +5.70%  '

# Row 3
$ws.Range("D3").Value = '1.705.82'
$ws.Range("E3").Value = '  +3.35%  '

# Row 4
$ws.Range("E4").Value = '  -0.37%  '

# Row 5
$ws.Range("D5").Value = '''331.25'
$ws.Range("E5").Value = '  +6.51%  '

# Row 6
$ws.Range("D6").Value = '''0.9972'
$ws.Range("E6").Value = '  -0.35%  '

# Row 7
$ws.Range("D7").Value = '''0.3679'
$ws.Range("E7").Value = '  +1.13%  '

# Row 8
$ws.Range("D8").Value = '''48.57'
$ws.Range("E8").Value = '  +4.02%  '

# Row 9
$ws.Range("D9").Value = '''0.3304'
$ws.Range("E9").Value = '  +1.88%  '

# Row 10
$ws.Range("D10").Value = '''1.168'
$ws.Range("E10").Value = '  +4.28%  '

# Row 11
$ws.Range("D11").Value = '''0.07335'
$ws.Range("E11").Value = '  +4.73%  '

# Row 12
$ws.Range("D12").Value = '''0.9986'
$ws.Range("E12").Value = '  -0.23%  '

# Row 13
$ws.Range("D13").Value = '''6.195'
$ws.Range("E13").Value = '  +4.43%  '

# Row 14
$ws.Range("D14").Value = '''19.91'
$ws.Range("E14").Value = '  +2.90%  '

# Row 15
$ws.Range("D15").Value = '''6.855'
$ws.Range("E15").Value = '  +4.36%  '

# Row 16
$ws.Range("D16").Value = '1.698.74'
$ws.Range("E16").Value = '  +3.26%  '

# Row 17
$ws.Range("D17").Value = '''0.00001070'
$ws.Range("E17").Value = '  +3.48%  '

# Row 18
$ws.Range("D18").Value = '''0.06615'
$ws.Range("E18").Value = '  +0.15%  '

# Row 19
$ws.Range("D19").Value = '''81.12'
$ws.Range("E19").Value = '  +3.84%  '

# Row 20
$ws.Range("E20").Value = '  -0.36%  '

# Row 21
$ws.Range("D21").Value = '''16.16'
$ws.Range("E21").Value = '  +4.00%  '

# Row 22
$ws.Range("D22").Value = '''6.045'
$ws.Range("E22").Value = '  +2.26%  '

# Row 23
$ws.Range("D23").Value = '''12.95'
$ws.Range("E23").Value = '  +3.85%  '

# Row 24
$ws.Range("D24").Value = '25.753.66'
$ws.Range("E24").Value = '  +5.66%  '

# Row 25
$ws.Range("D25").Value = '''2.462'
$ws.Range("E25").Value = '  -0.55%  '

# Row 26
$ws.Range("D26").Value = '''2.486'
$ws.Range("E26").Value = '  +7.95%  '

# Row 27
$ws.Range("D27").Value = '''149.58'
$ws.Range("E27").Value = '  +2.29%  '

# Row 28
$ws.Range("D28").Value = '''19.14'
$ws.Range("E28").Value = '  +3.46%  '

# Row 29
$ws.Range("D29").Value = '''1.292'
$ws.Range("E29").Value = '  +9.52%  '

# Row 30
$ws.Range("D30").Value = '1.890.36'
$ws.Range("E30").Value = '  +3.32%  '

# Row 31
$ws.Range("D31").Value = '''128.11'
$ws.Range("E31").Value = '  +3.62%  '

# Row 32
$ws.Range("D32").Value = '''4.116'
$ws.Range("E32").Value = '  +1.01%  '

# Row 33
$ws.Range("D33").Value = '''5.957'
$ws.Range("E33").Value = '  +6.08%  '

# Row 34
$ws.Range("D34").Value = '''0.08472'
$ws.Range("E34").Value = '  +0.51%  '

# Row 35
$ws.Range("D35").Value = '''1.673'
$ws.Range("E35").Value = '  +0.59%  '

# Row 36
$ws.Range("D36").Value = '''12.80'
$ws.Range("E36").Value = '  +5.67%  '

# Row 37
$ws.Range("D37").Value = '''5.309'
$ws.Range("E37").Value = '  +3.17%  '

# Row 38
$ws.Range("D38").Value = '''1.274'
$ws.Range("E38").Value = '  +2.92%  '

# Row 39
$ws.Range("D39").Value = '''0.06221'
$ws.Range("E39").Value = '  +4.02%  '

# Row 40
$ws.Range("D40").Value = '''8.529'
$ws.Range("E40").Value = '  +5.45%  '

# Row 41
$ws.Range("E41").Value = '  +3.40%  '

# Row 42
$ws.Range("D42").Value = '''0.02254'
$ws.Range("E42").Value = '  +1.90%  '

# Row 43
$ws.Range("D43").Value = '''14.51'
$ws.Range("E43").Value = '  +16.32%  '

# Row 44
$ws.Range("D44").Value = '''0.6104'
$ws.Range("E44").Value = '  +4.08%  '

# Row 45
$ws.Range("D45").Value = '''0.9975'
$ws.Range("E45").Value = '  -0.29%  '

# Row 46
$ws.Range("D46").Value = '''3.844'
$ws.Range("E46").Value = '  +2.31%  '

# Row 47
$ws.Range("D47").Value = '''0.5834'
$ws.Range("E47").Value = '  +4.61%  '

# Row 48
$ws.Range("D48").Value = '''126.14'
$ws.Range("E48").Value = '  +3.30%  '

# Row 49
$ws.Range("D49").Value = '''2.003'
$ws.Range("E49").Value = '  +3.32%  '

# Row 50
$ws.Range("E50").Value = '  +4.89%  '

# Row 51
$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
$ws.Range("D51").Value = '''1.208'
$ws.Range("E51").Value = '  +2.87%  '
